$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Correct the spelling of the shared string used across both sheets
# ("HeroBanner" -> "Herobanner"). It appears in sheet1 (C4, E4) and sheet2 (A4).
$ws1.Range("C4").Value = "Herobanner"
$ws1.Range("E4").Value = "Herobanner"
$ws2.Range("A4").Value = "Herobanner"

# Update selections / active view state:
# Testdata sheet selection moves to B18 (no longer the selected tab)
$ws2.Range("B18").Select()

# TC01_Verify_HomePage becomes the active/selected sheet with A20 selected
$ws1.Activate()
$ws1.Range("A20").Select()
